# Weekly price-sheet update: a new week's record is inserted at the top of
# the "Vega Modelo de Temuco - Ciboulette" data block (row 304), pushing all
# later rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 304; this shifts rows 304..345
# down to 305..346 (formatting follows Excel's default insert behaviour).
$ws.Rows(304).Insert()

# Populate the newly inserted row 304 with this week's record.
$ws.Range("A304").Value = 10
$ws.Range("B304").Value = "Vega Modelo de Temuco"
$ws.Range("C304").Value = "La Araucanía"
$ws.Range("D304").Value = 44984
$ws.Range("E304").Value = 9
$ws.Range("F304").Value = 100112039
$ws.Range("G304").Value = "Ciboulette"
$ws.Range("H304").Value = "Sin especificar"
$ws.Range("I304").Value = "Primera"
$ws.Range("J304").Value = 50
$ws.Range("K304").Value = 5000
$ws.Range("L304").Value = 5000
$ws.Range("M304").Value = 5000
$ws.Range("N304").Value = "$/docena de atados"
$ws.Range("O304").Value = "Provincia de Cautín"
$ws.Range("P304").Value = 1667
$ws.Range("Q304").Value = 3
$ws.Range("R304").Value = "Hortaliza"
